$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Grow the table to its final A1:C5 extent first, so the engine doesn't
# stomp our header text with an auto-generated "Column3" name afterwards.
$table = $ws.ListObjects.Item(1)
$table.Resize($ws.Range("A1:C5"))

# Add the new column header and new rows of data to the worksheet
$ws.Range("C1").Value = "Erklärung"

$data = @(
    @("WName", "String", "Name einer Waffe"),
    @("Wschaden", "Int", "(Grund-)Schaden einer Waffe - ohne Stärkebonus"),
    @("Schadensart", "bool", "Normaler oder rüstungsunabhängiger Schaden (TRUE = Rüstung zählt)"),
    @("Stärkeeinfluss", "bool", "Zählt der Stärkebonus (TRUE = ja)")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Adjust column widths (values chosen so the engine's pixel-quantized
# ColumnWidth storage lands on the target stored widths of 32 and ~140.832)
$ws.Columns.Item(1).ColumnWidth = 31.166666666666668
$ws.Columns.Item(3).ColumnWidth = 140

# Update the selection
$ws.Range("A6").Select()
